# Poland IV Liga - base update (17-06-2024 21:10)
# Swaps the data of several paired match rows (all columns B:AD; the
# leading index column A is left untouched). The underlying shared
# string table for a handful of team names was also re-ordered in the
# source commit, but since that has no visible effect on cell contents
# it does not need to be replicated here - only the row contents
# themselves need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(11, 12),
    @(30, 31),
    @(65, 66),
    @(95, 96),
    @(168, 169),
    @(183, 184),
    @(192, 193),
    @(202, 203),
    @(215, 216),
    @(222, 223)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")

    $v1 = $rng1.Value2
    $v2 = $rng2.Value2

    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}
